$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11
$ws.Range("C2").Value = 1.83
$ws.Range("D2").Value = 6.3
$ws.Range("E2").Value = 6.3
$ws.Range("F2").Value = 11.5
$ws.Range("G2").Value = 1.05
$ws.Range("H2").Value = 1.92
$ws.Range("I2").Value = 1.05
$ws.Range("J2").Value = 1.92

$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 4.17
$ws.Range("D3").Value = 16.4
$ws.Range("E3").Value = 13.4
$ws.Range("F3").Value = 24.4
$ws.Range("G3").Value = 2.73
$ws.Range("H3").Value = 4.57
$ws.Range("I3").Value = 2.23
$ws.Range("J3").Value = 4.07

$ws.Range("B4").Value = 12
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 10.1
$ws.Range("E4").Value = 9.300000000000001
$ws.Range("F4").Value = 17.5
$ws.Range("G4").Value = 1.68
$ws.Range("H4").Value = 3.04
$ws.Range("I4").Value = 1.55
$ws.Range("J4").Value = 2.91

$ws.Range("C5").Value = 2.83
$ws.Range("D5").Value = 10.1
$ws.Range("E5").Value = 9.300000000000001
$ws.Range("F5").Value = 17.6
$ws.Range("G5").Value = 1.68
$ws.Range("H5").Value = 3.07
$ws.Range("I5").Value = 1.55
$ws.Range("J5").Value = 2.94

$ws.Range("C6").Value = 1.33
$ws.Range("D6").Value = 9.199999999999999
$ws.Range("E6").Value = 6.8
$ws.Range("F6").Value = 12.5
$ws.Range("G6").Value = 1.53
$ws.Range("H6").Value = 2.47
$ws.Range("I6").Value = 1.13
$ws.Range("J6").Value = 2.08

$ws.Range("C7").Value = 0.83
$ws.Range("D7").Value = 7.4
$ws.Range("E7").Value = 7.4
$ws.Range("F7").Value = 13
$ws.Range("G7").Value = 1.23
$ws.Range("H7").Value = 2.17
$ws.Range("I7").Value = 1.23
$ws.Range("J7").Value = 2.17

$ws.Range("B8").Value = 6
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 6.4
$ws.Range("E8").Value = 6.4
$ws.Range("F8").Value = 11.4
$ws.Range("G8").Value = 1.07
$ws.Range("H8").Value = 1.9
$ws.Range("I8").Value = 1.07
$ws.Range("J8").Value = 1.9
$ws.Range("L8").Value = 1.4

$ws.Range("C9").Value = 0.67
$ws.Range("D9").Value = 6.8
$ws.Range("E9").Value = 6.8
$ws.Range("F9").Value = 12.6
$ws.Range("G9").Value = 1.13
$ws.Range("H9").Value = 2.09
$ws.Range("I9").Value = 1.13
$ws.Range("J9").Value = 2.09

$ws.Range("C10").Value = 1.5
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = 7.4
$ws.Range("F10").Value = 13.6
$ws.Range("G10").Value = 1.5
$ws.Range("H10").Value = 2.54
$ws.Range("I10").Value = 1.23
$ws.Range("J10").Value = 2.27

$ws.Range("B11").Value = 11
$ws.Range("C11").Value = 1.83
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = 9
$ws.Range("F11").Value = 16.6
$ws.Range("G11").Value = 1.51
$ws.Range("H11").Value = 2.76
$ws.Range("I11").Value = 1.51
$ws.Range("J11").Value = 2.76

$ws.Range("B12").Value = 12
$ws.Range("D12").Value = 10.1
$ws.Range("E12").Value = 8.5
$ws.Range("F12").Value = 13.8
$ws.Range("G12").Value = 1.68
$ws.Range("H12").Value = 2.57
$ws.Range("I12").Value = 1.42
$ws.Range("J12").Value = 2.3

$ws.Range("C13").Value = 0.83
$ws.Range("D13").Value = 8.1
$ws.Range("E13").Value = 7.3
$ws.Range("F13").Value = 13
$ws.Range("G13").Value = 1.35
$ws.Range("H13").Value = 2.3
$ws.Range("I13").Value = 1.22
$ws.Range("J13").Value = 2.16

$ws.Range("B14").Value = 8
$ws.Range("C14").Value = 1.33
$ws.Range("D14").Value = 11.1
$ws.Range("E14").Value = 10.4
$ws.Range("F14").Value = 16.7
$ws.Range("G14").Value = 1.86
$ws.Range("H14").Value = 2.92
$ws.Range("I14").Value = 1.73
$ws.Range("J14").Value = 2.78
$ws.Range("L14").Value = 2.6

$ws.Range("C15").Value = 1.33
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = 5.4
$ws.Range("F15").Value = 9.5
$ws.Range("G15").Value = 1.16
$ws.Range("H15").Value = 1.85
$ws.Range("I15").Value = 0.9
$ws.Range("J15").Value = 1.58

$ws.Range("B16").Value = 8
$ws.Range("C16").Value = 1.33
$ws.Range("D16").Value = 7.9
$ws.Range("E16").Value = 6.3
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 1.31
$ws.Range("H16").Value = 2.26
$ws.Range("I16").Value = 1.05
$ws.Range("J16").Value = 1.99

$ws.Range("C17").Value = 1.33
$ws.Range("D17").Value = 7.9
$ws.Range("E17").Value = 7.1
$ws.Range("F17").Value = 11.9
$ws.Range("G17").Value = 1.32
$ws.Range("H17").Value = 2.11
$ws.Range("I17").Value = 1.18
$ws.Range("J17").Value = 1.98

$ws.Range("B18").Value = 9
$ws.Range("C18").Value = 1.5
$ws.Range("D18").Value = 8.699999999999999
$ws.Range("E18").Value = 6.3
$ws.Range("F18").Value = 11.6
$ws.Range("G18").Value = 1.44
$ws.Range("H18").Value = 2.32
$ws.Range("I18").Value = 1.06
$ws.Range("J18").Value = 1.93

$ws.Range("B19").Value = 8
$ws.Range("C19").Value = 1.33
$ws.Range("D19").Value = 9.5
$ws.Range("E19").Value = 8.699999999999999
$ws.Range("F19").Value = 16.2
$ws.Range("G19").Value = 1.58
$ws.Range("H19").Value = 2.83
$ws.Range("I19").Value = 1.45
$ws.Range("J19").Value = 2.69
